$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Date - 02.11.2023"
$ws.Range("C4").Value = "DIT.E35KT6.SLP3.ZIP : last exported on 01.11.23"
$ws.Range("D4").Value = "DEL.N1LN3L.WD.D.GPD.DMP : last exported on 02.11.23`nDEL.N1LN3L.WD.D.GPL.DMP : last exported on 02.11.23`nDEL.N1LN3L.WD.D.GTR.DMP : last exported on 02.11.23`nDEL.N1LKQQ.WD.D.UPD.DMP : last exported on 02.11.23`nDEL.N1LN3L.WD.D.ASPD.DMP : last exported on 02.11.23`nDEL.N1LN3L.WD.D.ASPL.DMP : last exported on 02.11.23`nDEL.N1LK61.WD.D.CPD.DMP : last exported on 02.11.23`nDEL.N1LN0R.WD.D.CPL.DMP : last exported on 02.11.23`nDEL.KT6KTZ.SL.D.GGO.ZIP : last exported on 02.11.23`nDEL.KT6KTZ.SL.D.UGO.ZIP : last exported on 02.11.23`nDEL.KT6KTZ.SL.D.ASGO.ZIP : last exported on 02.11.23`nDEL.KT6N5M.SL.D.GGO.ZIP : last exported on 02.11.23`nDEL.KT6N5M.SL.D.UGO.ZIP : last exported on 02.11.23`nDEL.KT6N5M.SL.D.ASGO.ZIP : last exported on 02.11.23"
$ws.Range("C7").Value = "DLV.R31KT6.WI.C.ZIP : last exported on 31.10.23`nDLV.KZ6KT6.WI.V.ZIP : last exported on 01.11.23`nDLV.RPKKT6.WI.S.ZIP : last exported on 31.10.23`nDLV.I5XKT6.WI.A.ZIP : last exported on 01.11.23`nDIT.E35KT6.WI.ZIP : last exported on 01.11.23`nDEL.K2PKT6.WI.ZIP : last exported on 01.11.23`nDEL.R7AKT6.WI.ZIP : last exported on 01.11.23`nDEL.N5FKT6.WI.ZIP : last exported on 27.10.23"
$ws.Range("D7").Value = "DEL.N1LN3L.WI.D.GPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.WI.D.GPL.DMP : last exported on 01.11.23`nDEL.N1LN3L.WI.D.GTR.DMP : last exported on 01.11.23`nDEL.N1LN3L.WI.D.ASPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.WI.D.ASPL.DMP : last exported on 01.11.23`nDEL.N1LK61.WI.D.CPD.DMP : last exported on 01.11.23`nDEL.N1LN0R.WI.D.CPL.DMP : last exported on 01.11.23`nDEL.KT6KTZ.WI.GGO.ZIP : last exported on 01.11.23`nDEL.KT6KTZ.WI.IMG.GGO.ZIP : last exported on 01.11.23`nDEL.KT6E35.WI.GGO.ZIP : last exported on 01.11.23`nDEL.KT6E35.WI.IMG.GGO.ZIP : last exported on 01.11.23`nDEL.KT6N5M.WI.GGO.ZIP : last exported on 01.11.23`nDEL.KT6N5M.WI.IMG.GGO.ZIP : last exported on 01.11.23"
$ws.Range("C10").Value = "DEL.N3LKT6.HST.MASTER.ZIP : last exported on 21.10.23`nDEL.N3LKT6.HST.??????.ZIP : last exported on 02.11.23`nDGV.N3LKT6.EPELS.??????.ZIP : last exported on 30.10.23"
$ws.Range("D10").Value = "DEL.N1LN3L.TP.D.GPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.TP.D.GPL.DMP : last exported on 01.11.23`nDEL.N1LN3L.TP.D.GTR.DMP : last exported on 01.11.23`nDEL.N1LKQQ.TP.D.UPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.TP.D.ASPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.TP.D.ASPL.DMP : last exported on 01.11.23`nDEL.N1LK61.TP.D.CPD.DMP : last exported on 01.11.23`nDEL.N1LN0R.TP.D.CPL.DMP : last exported on 01.11.23`nDEL.KT6E35.TP.D.GGO.ZIP : last exported on 01.11.23`nDEL.KT6E35.TP.D.UGO.ZIP : last exported on 01.11.23`nDEL.KT6N5M.TP.D.GGO.ZIP : last exported on 01.11.23`nDEL.KT6N5M.TP.D.UGO.ZIP : last exported on 01.11.23`nDEL.KT6N5M.TP.D.ASGO.ZIP : last exported on 01.11.23`nDEL.KT6KTZ.TP.D.GGO.ZIP : last exported on 01.11.23`nDEL.KT6KTZ.TP.D.UGO.ZIP : last exported on 01.11.23`nDEL.KT6KTZ.TP.D.ASGO.ZIP : last exported on 01.11.23`nDED.KT6I1M.HST2.D.GGO.ZIP : last exported on 01.11.23"
$ws.Range("D12").Value = "DEL.N1LN3L.MC.F.GPD.DMP : last exported on 28.10.23`nDEL.N1LN3L.MC.F.GPL.DMP : last exported on 28.10.23`nDEL.N1LN3L.MC.F.GTR.DMP : last exported on 01.11.23`nDEL.N1LKQQ.MC.F.UPD.DMP : last exported on 28.10.23`nDEL.N1LN3L.MC.F.ASPD.DMP : last exported on 28.10.23`nDEL.N1LN3L.MC.F.ASPL.DMP : last exported on 28.10.23`nDGI.KT6R11.MANDANT.TXT : last exported on 28.10.23"
$ws.Range("C13").Value = "DDC.R11KT6.ELFI.TT.TXT : last exported on 01.11.23"
$ws.Range("D14").Value = "DEL.N1LN3L.AU.D.GPD.DMP : last exported on 27.10.23`nDEL.N1LN3L.AU.D.GPL.DMP : last exported on 27.10.23`nDEL.N1LN3L.AU.D.GTR.DMP : last exported on 01.11.23`nDEL.KT6N5M.AU.GGO.ZIP : last exported on 27.10.23`nDEL.KT6N5M.AU.IMG.GGO.ZIP : last exported on 27.10.23`nDEL.KT6KTZ.AU.GGO.ZIP : last exported on 27.10.23`nDEL.KT6KTZ.AU.IMG.GGO.ZIP : last exported on 27.10.23`nDEL.KT6E35.AU.GGO.ZIP : last exported on 27.10.23`nDEL.KT6E35.AU.IMG.GGO.ZIP : last exported on 27.10.23"
$ws.Range("C15").Value = "DGQ.R11KT6.BSIVM.TXT : last exported on 01.11.23`nDGQ.R11KT6.BSIVM.C.TXT : last exported on 01.11.23`nDGQ.R11KT6.BSIVM.US2.TXT : last exported on 01.11.23`nDGQ.R11KT6.PKAT.TXT : last exported on 01.11.23`nDGQ.R11KT6.PKAT.C.TXT : last exported on 01.11.23`nDGQ.R11KT6.PKAT.US2.TXT : last exported on 01.11.23`nDGQ.R11KT6.BSITK.TXT : last exported on 01.11.23`nDGQ.R11KT6.BSITX.TXT : last exported on 01.11.23`nDGQ.R31KT6.BSIVMC.TXT : last exported on 10.10.23`nDGQ.R31KT6.PKATC.TXT : last exported on 10.10.23`nDEL.KMPKT6.APOS.DATA.ZIP : last exported on 01.11.23"
$ws.Range("D15").Value = "DEL.N1LN3L.AP.F.GPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.AP.F.GPL.DMP : last exported on 01.11.23`nDEL.N1LN3L.AP.F.GTR.DMP : last exported on 01.11.23`nDEL.N1LKQQ.AP.F.UPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.AP.F.ASPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.AP.F.ASPL.DMP : last exported on 01.11.23`nDEL.N1LK61.AP.F.CPD.DMP : last exported on 01.11.23`nDEL.N1LN0R.AP.F.CPL.DMP : last exported on 01.11.23`nDEL.N1LR31.AP.F.PPS.DMP : last exported on 01.11.23`nDEL.KT6E35.AP.F.DB.ZIP : last exported on 01.11.23`nDEL.KT6KGZ.AP.F.DWH.ZIP : last exported on 01.11.23`nDEL.KT6UGV.AP.F.CEB.ZIP : last exported on 01.11.23`nDVL.KT6N5X.VLM.AP.ZIP : last exported on 01.11.23"
$ws.Range("D16").Value = "DGI.N1LR11.FB.ORGUNIT.CSV : last exported on 01.11.23`nDGI.N1LR11.FB.IMPTR.CSV : last exported on 01.11.23"
$ws.Range("C17").Value = "DHP.KUDKT6.KVPSIMPORT : last exported on 01.11.23`nDHP.KUDKT6.ORGUNITS : last exported on 02.11.23`nDHS.R11KT6.HSB02ALL : last exported on 01.11.23"
$ws.Range("D17").Value = "DEL.N1LN3L.DP.F.GPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.DP.F.GPL.DMP : last exported on 01.11.23`nDEL.N1LN3L.DP.F.GTR.DMP : last exported on 01.11.23`nDEL.N1LKQQ.DP.F.UPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.DP.F.ASPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.DP.F.ASPL.DMP : last exported on 01.11.23`nDEL.N1LK61.DP.F.CPD.DMP : last exported on 01.11.23`nDEL.N1LN0R.DP.F.CPL.DMP : last exported on 01.11.23`nDEL.N1LN3L.DISSORG.ZIP : last exported on 01.11.23`nDEL.N1LN3L.DISSAORG.ZIP : last exported on 01.11.23"
$ws.Range("C20").Value = "DEL.N3LKT6.WT.??????.ZIP : last exported on 01.11.23"
$ws.Range("D20").Value = "DEL.N1LN3L.MT.F.GPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.MT.F.GPL.DMP : last exported on 01.11.23`nDEL.N1LN3L.MT.F.GTR.DMP : last exported on 01.11.23`nDEL.N1LKQQ.MT.F.UPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.MT.F.ASPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.MT.F.ASPL.DMP : last exported on 01.11.23`nDEL.N1LK61.MT.F.CPD.DMP : last exported on 01.11.23`nDEL.N1LN0R.MT.F.CPL.DMP : last exported on 01.11.23`nDEL.N1LR31.MT.F.GPD.DMP : last exported on 01.11.23`nDEL.N1LI4X.MT.F.GPD.ZIP : last exported on 01.11.23`nDEL.N1LI4Y.MT.F.GPL.ZIP : last exported on 01.11.23`nDEL.N1LGN1.NL.F.GPD.ZIP : last exported on 01.11.23`nDEL.N1LKT6.MT.??????.G : last exported on 01.11.23`nDEL.N1LKT6.MT.??????.U : last exported on 01.11.23`nDEL.N1LKT6.MT.??????.AS : last exported on 01.11.23`nDVL.N1LN5X.VLM.DSP : last exported on 01.11.23`nDVL.N1LN5X.VLM.WT : last exported on 01.11.23`nDDS.N1LR11.DSP : last exported on 01.11.23`nDDS.N1LR11.WT : last exported on 01.11.23`nDEL.KT6KTZ.MT.F.GGO.ZIP : last exported on 01.11.23`nDEL.KT6KTZ.MT.F.UGO.ZIP : last exported on 01.11.23`nDEL.KT6KTZ.MT.F.ASGO.ZIP : last exported on 01.11.23`nDEL.KT6N5M.MT.F.GGO.ZIP : last exported on 01.11.23`nDEL.KT6N5M.MT.F.UGO.ZIP : last exported on 01.11.23`nDEL.KT6N5M.MT.F.ASGO.ZIP : last exported on 01.11.23`nDEL.KT6E35.MT.F.GGO.ZIP : last exported on 01.11.23`nDEL.KT6E35.MT.F.UGO.ZIP : last exported on 01.11.23`nDEL.KT6N5X.MT.F.GGO.ZIP : last exported on 01.11.23`nDEL.KT6N5X.MT.F.UGO.ZIP : last exported on 01.11.23`nDEL.KT6N5X.MT.F.ASGO.ZIP : last exported on 01.11.23"
$ws.Range("C21").Value = "DHP.KUDKT6.EXEMN.ZIP : last exported on 02.11.23`nDHP.KUDKT6.DPO.EMAILS : last exported on 02.11.23"
$ws.Range("D21").Value = "DEL.N1LN3L.EM.F.GPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.EM.F.GPL.DMP : last exported on 01.11.23`nDEL.N1LN3L.EM.F.GTR.DMP : last exported on 01.11.23`nDEL.N1LKQQ.EM.F.UPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.EM.F.ASPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.EM.F.ASPL.DMP : last exported on 01.11.23`nDEL.N1LK61.EM.F.CPD.DMP : last exported on 01.11.23`nDEL.N1LN0R.EM.F.CPL.DMP : last exported on 01.11.23`nDEL.N1LN3L.EM.F.DISS.ZIP : last exported on 01.11.23`nDEL.N1LN3L.EM.F.DISSA.ZIP : last exported on 01.11.23"
$ws.Range("C22").Value = "DKG.R11KT6.L520.P.EDCP.HVB : last exported on 02.11.23`nDKG.R11KT6.L520.P.EDCP.LEI : last exported on 02.11.23`nDKG.R11KT6.L520.P.EDCP.SGE : last exported on 02.11.23`nDKG.R11KT6.L520.P.EDCP.SMO : last exported on 02.11.23`nDKG.R11KT6.L520.P.EDCP.STR : last exported on 02.11.23`nDKG.R11KT6.L520.P.EDCP.TEC : last exported on 02.11.23`nDKG.R11KT6.L520.P.EDCP.TGE : last exported on 02.11.23`nDKG.R11KT6.L520.P.EDCP.TMO : last exported on 02.11.23`nDDC.R11KT6.ELFI.MD.TXT : last exported on 29.10.23`nDDC.R11KT6.ELFI.TT.EXT.TXT : last exported on 01.11.23`nDDC.R11KT6.ELFI.PK.TXT : last exported on 29.10.23`nDDC.R11KT6.ELFI.PR.TXT : last exported on 29.10.23"
$ws.Range("D22").Value = "DEL.N1LN3L.EL.F.GPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.EL.F.GPL.DMP : last exported on 01.11.23`nDEL.N1LE35.ELFI.G.DATA : last exported on 01.11.23`nDEL.N1LN5M.EL.F.GGO.ZIP : last exported on 01.11.23`nDEL.N1LKTZ.EL.F.GGO.ZIP : last exported on 01.11.23`nDEL.N1LN5X.EL.F.GGO.ZIP : last exported on 01.11.23`nDEL.N1LE35.EL.F.GGO.ZIP : last exported on 01.11.23`nDEL.N1LKT6.EC.??????.G : last exported on 01.11.23`nDWT.N1LN3L.WT.??????.ZIP : last exported on 01.11.23`nDEL.N1LKQQ.EL.F.UPD.DMP : last exported on 01.11.23`nDEL.N1LE35.ELFI.U.DATA : last exported on 01.11.23`nDEL.N1LN5M.EL.F.UGO.ZIP : last exported on 01.11.23`nDEL.N1LKTZ.EL.F.UGO.ZIP : last exported on 01.11.23`nDEL.N1LN5X.EL.F.UGO.ZIP : last exported on 01.11.23`nDEL.N1LE35.EL.F.UGO.ZIP : last exported on 01.11.23`nDEL.N1LKT6.EC.??????.U : last exported on 01.11.23`nDEL.N1LN3L.EL.F.ASPD.DMP : last exported on 01.11.23`nDEL.N1LN3L.EL.F.ASPL.DMP : last exported on 01.11.23`nDEL.N1LK61.EL.F.CPD.DMP : last exported on 01.11.23`nDEL.N1LN0R.EL.F.CPL.DMP : last exported on 01.11.23`nDEL.N1LN5M.EL.F.ASGO.ZIP : last exported on 01.11.23`nDEL.N1LKTZ.EL.F.ASGO.ZIP : last exported on 01.11.23`nDEL.N1LN5X.EL.F.ASGO.ZIP : last exported on 01.11.23`nDEL.N1LE35.EL.F.ASGO.ZIP : last exported on 01.11.23`nDEL.N1LKT6.EC.??????.AS : last exported on 01.11.23`nDEL.N1LN3L.EL.F.GTR.DMP : last exported on 01.11.23`nDGV.N1LR31.EL.F.HST.ZIP : last exported on 01.11.23`nDGV.N1LN3L.EL.F.HST.ZIP : last exported on 01.11.23`nDGV.N1LCEN.EL.F.ZIP : last exported on 01.11.23`nDGV.N1LN7K.EL.F.ZIP : last exported on 01.11.23`nDGV.N1LN5X.EL.F.ZIP : last exported on 01.11.23`nDEL.KT6KT6.EL.F.DRE.CSV : last exported on 01.11.23`nDEL.KT6N5X.EL.F.AU.ZIP : last exported on 01.11.23"
